$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 28 (hunk 0)
$ws.Range("H28").Value = 2299.8
$ws.Range("I28").Value = 762.4167
$ws.Range("J28").Value = 4605.875
$ws.Range("K28").Value = 762.4167
$ws.Range("L28").Value = 4605.875
$ws.Range("M28").Value = -277.4167
$ws.Range("N28").Value = -5575.875
# row 64 (hunk 1)
$ws.Range("H64").Value = 7894.5
$ws.Range("I64").Value = 7855.952
$ws.Range("J64").Value = 8164.3335
$ws.Range("K64").Value = 7855.952
$ws.Range("L64").Value = 8164.3335
$ws.Range("M64").Value = -7607.952
$ws.Range("N64").Value = -8660.333500000001
# row 67 (hunk 2)
$ws.Range("H67").Value = 7894.5
$ws.Range("I67").Value = 7855.952
$ws.Range("J67").Value = 8164.3335
$ws.Range("K67").Value = 7855.952
$ws.Range("L67").Value = 8164.3335
$ws.Range("M67").Value = -6997.952
$ws.Range("N67").Value = -9880.333500000001
# row 74 (hunk 3)
$ws.Range("H74").Value = 5023.6875
$ws.Range("I74").Value = 4789.5
$ws.Range("J74").Value = 5726.25
$ws.Range("K74").Value = 4789.5
$ws.Range("L74").Value = 5726.25
$ws.Range("M74").Value = -3853.5
$ws.Range("N74").Value = -7598.25
# row 77 (hunk 4)
$ws.Range("H77").Value = 5023.6875
$ws.Range("I77").Value = 4789.5
$ws.Range("J77").Value = 5726.25
$ws.Range("K77").Value = 23947.5
$ws.Range("L77").Value = 28631.25
$ws.Range("M77").Value = -19267.5
$ws.Range("N77").Value = -37991.25
# row 86 (hunk 5)
$ws.Range("H86").Value = 2683.7666
$ws.Range("I86").Value = 2224.5
$ws.Range("J86").Value = 3372.6667
$ws.Range("K86").Value = 2224.5
$ws.Range("L86").Value = 3372.6667
$ws.Range("M86").Value = -1101.5
$ws.Range("N86").Value = -5618.6667
# row 89 (hunk 6)
$ws.Range("H89").Value = 2683.7666
$ws.Range("I89").Value = 2224.5
$ws.Range("J89").Value = 3372.6667
$ws.Range("K89").Value = 11122.5
$ws.Range("L89").Value = 16863.3335
$ws.Range("M89").Value = -5506.5
$ws.Range("N89").Value = -28095.3335
# row 98 (hunk 7)
$ws.Range("H98").Value = 614.619
$ws.Range("I98").Value = 626.4211
$ws.Range("K98").Value = 626.4211
$ws.Range("M98").Value = 871.5789
# row 122 (hunk 8)
$ws.Range("H122").Value = 614.619
$ws.Range("I122").Value = 626.4211
$ws.Range("K122").Value = 1879.2633
$ws.Range("M122").Value = 570.7366999999999
# row 138 (hunk 9)
$ws.Range("H138").Value = 1748.0454
$ws.Range("I138").Value = 1629.579
$ws.Range("J138").Value = 2498.3333
$ws.Range("K138").Value = 4888.737
$ws.Range("L138").Value = 7494.999899999999
$ws.Range("M138").Value = 251.2629999999999
$ws.Range("N138").Value = -17774.9999

$ws = $wb.Worksheets.Item("ARM")
# row 97 (hunk 10)
$ws.Range("H97").Value = 1998.3334
$ws.Range("I97").Value = 1998
$ws.Range("K97").Value = 1998
$ws.Range("M97").Value = -1502
# row 132 (hunk 11)
$ws.Range("H132").Value = 1805.0294
$ws.Range("I132").Value = 1705.5483
$ws.Range("K132").Value = 5116.644899999999
$ws.Range("M132").Value = -2586.644899999999

$ws = $wb.Worksheets.Item("BSM")
# row 99 (hunk 12)
$ws.Range("H99").Value = 1699479.6
$ws.Range("I99").Value = 112846
$ws.Range("J99").Value = 3127450
$ws.Range("K99").Value = 112846
$ws.Range("L99").Value = 3127450
$ws.Range("M99").Value = -111348
$ws.Range("N99").Value = -3130446
# row 105 (hunk 13)
$ws.Range("H105").Value = 25179.023
$ws.Range("I105").Value = 28471.676
$ws.Range("J105").Value = 4874.3335
$ws.Range("K105").Value = 28471.676
$ws.Range("L105").Value = 4874.3335
$ws.Range("M105").Value = -26724.676
$ws.Range("N105").Value = -8368.333500000001
# row 134 (hunk 14)
$ws.Range("H134").Value = 4619.5
$ws.Range("J134").Value = 6650.8
$ws.Range("L134").Value = 19952.4
$ws.Range("N134").Value = -25022.4

$ws = $wb.Worksheets.Item("CRP")
# row 31 (hunk 15)
$ws.Range("H31").Value = 3006.6956
$ws.Range("I31").Value = 1877.7142
$ws.Range("K31").Value = 1877.7142
$ws.Range("M31").Value = -1582.7142
# row 34 (hunk 16)
$ws.Range("H34").Value = 3006.6956
$ws.Range("I34").Value = 1877.7142
$ws.Range("K34").Value = 1877.7142
$ws.Range("M34").Value = -1675.7142
# row 99 (hunk 17)
$ws.Range("H99").Value = 10170642
$ws.Range("I99").Value = 13890636
$ws.Range("J99").Value = 5210650
$ws.Range("K99").Value = 13890636
$ws.Range("L99").Value = 5210650
$ws.Range("M99").Value = -13889138
$ws.Range("N99").Value = -5213646
# row 126 (hunk 18)
$ws.Range("H126").Value = 10170642
$ws.Range("I126").Value = 13890636
$ws.Range("J126").Value = 5210650
$ws.Range("K126").Value = 41671908
$ws.Range("L126").Value = 15631950
$ws.Range("M126").Value = -41669438
$ws.Range("N126").Value = -15636890

$ws = $wb.Worksheets.Item("GSM")
# row 97 (hunk 19)
$ws.Range("H97").Value = 421.85715
$ws.Range("I97").Value = 283.45456
$ws.Range("J97").Value = 929.3333
$ws.Range("K97").Value = 283.45456
$ws.Range("L97").Value = 929.3333
$ws.Range("M97").Value = 212.54544
$ws.Range("N97").Value = -1921.3333
# row 102 (hunk 20)
$ws.Range("H102").Value = 1889.6
$ws.Range("I102").Value = 1889.6
$ws.Range("K102").Value = 1889.6
$ws.Range("M102").Value = -267.5999999999999
# row 113 (hunk 21)
$ws.Range("H113").Value = 5559141
$ws.Range("J113").Value = 6669967
$ws.Range("L113").Value = 6669967
$ws.Range("N113").Value = -6674307
# row 132 (hunk 22)
$ws.Range("H132").Value = 4203.8076
$ws.Range("I132").Value = 3434.1365
$ws.Range("J132").Value = 8437
$ws.Range("K132").Value = 10302.4095
$ws.Range("L132").Value = 25311
$ws.Range("M132").Value = -7772.4095
$ws.Range("N132").Value = -30371

$ws = $wb.Worksheets.Item("LTW")
# row 26 (hunk 23)
$ws.Range("H26").Value = 9
$ws.Range("I26").Value = 9
$ws.Range("K26").Value = 9
$ws.Range("M26").Value = 286
# row 122 (hunk 24)
$ws.Range("H122").Value = 46158016
$ws.Range("I122").Value = 58827364
$ws.Range("J122").Value = 22227026
$ws.Range("K122").Value = 176482092
$ws.Range("L122").Value = 66681078
$ws.Range("M122").Value = -176479642
$ws.Range("N122").Value = -66685978
# row 136 (hunk 25)
$ws.Range("H136").Value = 2117.6086
$ws.Range("J136").Value = 3068.1428
$ws.Range("L136").Value = 9204.428400000001
$ws.Range("N136").Value = -14304.4284

$ws = $wb.Worksheets.Item("WVR")
# row 62 (hunk 26)
$ws.Range("H62").Value = 6666.6665
$ws.Range("I62").Value = 6666.6665
$ws.Range("K62").Value = 6666.6665
$ws.Range("M62").Value = -6042.6665
# row 65 (hunk 27)
$ws.Range("H65").Value = 6666.6665
$ws.Range("I65").Value = 6666.6665
$ws.Range("K65").Value = 33333.3325
$ws.Range("M65").Value = -30213.3325
# row 122 (hunk 28)
$ws.Range("H122").Value = 2129.5
$ws.Range("I122").Value = 2126
$ws.Range("K122").Value = 6378
$ws.Range("M122").Value = -3928
